$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.045889371829667
$ws.Range("D2").Value = 1.04687021693114
$ws.Range("E2").Value = 1.04350055558008
$ws.Range("F2").Value = 1.054336992649198
$ws.Range("I2").Value = 1.045805903132512
$ws.Range("J2").Value = 1.050946695619129
$ws.Range("K2").Value = 1.049634443028699
$ws.Range("L2").Value = 1.046274245442233
$ws.Range("M2").Value = 1.057080481744324

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.047133903105042
$ws.Range("D3").Value = 1.047825673414588
$ws.Range("E3").Value = 1.04457458345226
$ws.Range("F3").Value = 1.055789571921186
$ws.Range("I3").Value = 1.046257205973744
$ws.Range("J3").Value = 1.051837617078412
$ws.Range("K3").Value = 1.050401256093648
$ws.Range("L3").Value = 1.047158633558718
$ws.Range("M3").Value = 1.058344650361362

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.047938346416962
$ws.Range("D4").Value = 1.048443102837044
$ws.Range("E4").Value = 1.045268922342316
$ws.Range("F4").Value = 1.05672901742829
$ws.Range("I4").Value = 1.046547473026928
$ws.Range("J4").Value = 1.052412749953195
$ws.Range("K4").Value = 1.050896002187085
$ws.Range("L4").Value = 1.047729670680783
$ws.Range("M4").Value = 1.059161666603921

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.048276333986484
$ws.Range("D5").Value = 1.048702477078752
$ws.Range("E5").Value = 1.045560674189358
$ws.Range("F5").Value = 1.057123852329814
$ws.Range("I5").Value = 1.046669082303923
$ws.Range("J5").Value = 1.052654214313695
$ws.Range("K5").Value = 1.05110365230463
$ws.Range("L5").Value = 1.0479694445691
$ws.Range("M5").Value = 1.059504908402245

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.048333071936278
$ws.Range("D6").Value = 1.048746015889867
$ws.Range("E6").Value = 1.045609651984626
$ws.Range("F6").Value = 1.057190140648861
$ws.Range("I6").Value = 1.046689476507696
$ws.Range("J6").Value = 1.052694738458495
$ws.Range("K6").Value = 1.051138497725822
$ws.Range("L6").Value = 1.048009686711805
$ws.Range("M6").Value = 1.059562526734064

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.047942863406498
$ws.Range("D7").Value = 1.048446569364601
$ws.Range("E7").Value = 1.045272821324263
$ws.Range("F7").Value = 1.056734293650131
$ws.Range("I7").Value = 1.04654909962064
$ws.Range("J7").Value = 1.052415977672115
$ws.Range("K7").Value = 1.050898778156112
$ws.Range("L7").Value = 1.047732875688418
$ws.Range("M7").Value = 1.059166253923498

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.046310145440746
$ws.Range("D8").Value = 1.047193287830271
$ws.Range("E8").Value = 1.043863659977354
$ws.Range("F8").Value = 1.054827998806527
$ws.Range("I8").Value = 1.045958787497203
$ws.Range("J8").Value = 1.051248068258152
$ws.Range("K8").Value = 1.0498938890488
$ws.Range("L8").Value = 1.046573382729362
$ws.Range("M8").Value = 1.057507919932786

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.043426408670658
$ws.Range("D9").Value = 1.044978515848589
$ws.Range("E9").Value = 1.041375614666168
$ws.Range("F9").Value = 1.051465053358857
$ws.Range("I9").Value = 1.044905069173747
$ws.Range("J9").Value = 1.049179593983745
$ws.Range("K9").Value = 1.048112080068681
$ws.Range("L9").Value = 1.044520755641526
$ws.Range("M9").Value = 1.054577993936813

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041499205895602
$ws.Range("D10").Value = 1.043497618111978
$ws.Range("E10").Value = 1.039713457489666
$ws.Range("F10").Value = 1.049220235897897
$ws.Range("I10").Value = 1.044193414691789
$ws.Range("J10").Value = 1.047793420712012
$ws.Range("K10").Value = 1.046916640737115
$ws.Range("L10").Value = 1.043145840574493
$ws.Range("M10").Value = 1.052619255702356

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040663538485794
$ws.Range("D11").Value = 1.042855304358263
$ws.Range("E11").Value = 1.038992870159661
$ws.Range("F11").Value = 1.048247462768925
$ws.Range("I11").Value = 1.043883063598791
$ws.Range("J11").Value = 1.047191452485614
$ws.Range("K11").Value = 1.046397178827657
$ws.Range("L11").Value = 1.042548913992351
$ws.Range("M11").Value = 1.051769748957502

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040352953355364
$ws.Range("D12").Value = 1.042616556284926
$ws.Range("E12").Value = 1.038725079351284
$ws.Range("F12").Value = 1.047886012761749
$ws.Range("I12").Value = 1.04376745302067
$ws.Range("J12").Value = 1.046967589303075
$ws.Range("K12").Value = 1.046203950224641
$ws.Range("L12").Value = 1.042326948751028
$ws.Range("M12").Value = 1.051453994916261

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040419583183463
$ws.Range("D13").Value = 1.042667776026999
$ws.Range("E13").Value = 1.038782527437657
$ws.Range("F13").Value = 1.0479635504813
$ws.Range("I13").Value = 1.043792266951896
$ws.Range("J13").Value = 1.047015620773926
$ws.Range("K13").Value = 1.046245411006618
$ws.Range("L13").Value = 1.042374571961278
$ws.Range("M13").Value = 1.051521734787835

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040637869140737
$ws.Range("D14").Value = 1.042835572731676
$ws.Range("E14").Value = 1.03897073721033
$ws.Range("F14").Value = 1.048217587654319
$ws.Range("I14").Value = 1.043873513983123
$ws.Range("J14").Value = 1.04717295331653
$ws.Range("K14").Value = 1.046381212163206
$ws.Range("L14").Value = 1.042530571191452
$ws.Range("M14").Value = 1.051743652931262

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.040772338201522
$ws.Range("D15").Value = 1.042938935990385
$ws.Range("E15").Value = 1.03908668180236
$ws.Range("F15").Value = 1.048374092433185
$ws.Range("I15").Value = 1.04392352886399
$ws.Range("J15").Value = 1.047269855831457
$ws.Range("K15").Value = 1.046464846913091
$ws.Range("L15").Value = 1.042626655560327
$ws.Range("M15").Value = 1.051880356122685

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041554641241166
$ws.Range("D16").Value = 1.043540223484309
$ws.Range("E16").Value = 1.039761262070842
$ws.Range("F16").Value = 1.049284779307709
$ws.Range("I16").Value = 1.044213965163463
$ws.Range("J16").Value = 1.047833334312869
$ws.Range("K16").Value = 1.046951076964055
$ws.Range("L16").Value = 1.043185423108246
$ws.Range("M16").Value = 1.052675605578558

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.0420450413628
$ws.Range("D17").Value = 1.043917105534409
$ws.Range("E17").Value = 1.040184175447482
$ws.Range("F17").Value = 1.049855823665824
$ws.Range("I17").Value = 1.044395557882023
$ws.Range("J17").Value = 1.048186319856692
$ws.Range("K17").Value = 1.047255584802616
$ws.Range("L17").Value = 1.043535498393639
$ws.Range("M17").Value = 1.053174076830344

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042330970298462
$ws.Range("D18").Value = 1.044136830888229
$ws.Range("E18").Value = 1.040430770672201
$ws.Range("F18").Value = 1.050188831976002
$ws.Range("I18").Value = 1.044501265698902
$ws.Range("J18").Value = 1.048392042091707
$ws.Range("K18").Value = 1.047433022796065
$ws.Range("L18").Value = 1.043739539201533
$ws.Range("M18").Value = 1.053464695786015

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042428445614691
$ws.Range("D19").Value = 1.044211734066494
$ws.Range("E19").Value = 1.040514839285501
$ws.Range("F19").Value = 1.050302367129555
$ws.Range("I19").Value = 1.04453727340474
$ws.Range("J19").Value = 1.048462159609453
$ws.Range("K19").Value = 1.047493494743892
$ws.Range("L19").Value = 1.0438090861172
$ws.Range("M19").Value = 1.053563767290508

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.041992437808439
$ws.Range("D20").Value = 1.04387668038974
$ws.Range("E20").Value = 1.040138809452025
$ws.Range("F20").Value = 1.049794563521773
$ws.Range("I20").Value = 1.044376096658848
$ws.Range("J20").Value = 1.048148465275615
$ws.Range("K20").Value = 1.04722293224601
$ws.Range("L20").Value = 1.043497954391347
$ws.Range("M20").Value = 1.053120609159386

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.040573594400016
$ws.Range("D21").Value = 1.042786165336576
$ws.Range("E21").Value = 1.038915317807897
$ws.Range("F21").Value = 1.04814278336286
$ws.Range("I21").Value = 1.043849597947129
$ws.Range("J21").Value = 1.04712663015457
$ws.Range("K21").Value = 1.046341229783433
$ws.Range("L21").Value = 1.042484639958172
$ws.Range("M21").Value = 1.051678309390758

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.039680461572222
$ws.Range("D22").Value = 1.042099564241416
$ws.Range("E22").Value = 1.038145291410309
$ws.Range("F22").Value = 1.047103553689462
$ws.Range("I22").Value = 1.043516643439253
$ws.Range("J22").Value = 1.046482624922181
$ws.Range("K22").Value = 1.045785262677194
$ws.Range("L22").Value = 1.041846138207117
$ws.Range("M22").Value = 1.050770264164667

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040154029146483
$ws.Range("D23").Value = 1.042463635488443
$ws.Range("E23").Value = 1.038553570849414
$ws.Range("F23").Value = 1.047654536315798
$ws.Range("I23").Value = 1.043693331899944
$ws.Range("J23").Value = 1.046824170988426
$ws.Range("K23").Value = 1.046080144427979
$ws.Range("L23").Value = 1.042184752851463
$ws.Range("M23").Value = 1.051251752899727

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.042016207439799
$ws.Range("D24").Value = 1.043894947092489
$ws.Range("E24").Value = 1.040159308648964
$ws.Range("F24").Value = 1.049822244566725
$ws.Range("I24").Value = 1.044384891003313
$ws.Range("J24").Value = 1.048165570651147
$ws.Range("K24").Value = 1.047237687075401
$ws.Range("L24").Value = 1.043514919380039
$ws.Range("M24").Value = 1.053144769302552

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.044172737943938
$ws.Range("D25").Value = 1.045551850659487
$ws.Range("E25").Value = 1.042019432846232
$ws.Range("F25").Value = 1.052334936899116
$ws.Range("I25").Value = 1.045179091234905
$ws.Range("J25").Value = 1.049715600343952
$ws.Range("K25").Value = 1.048574044759459
$ws.Range("L25").Value = 1.045052543649962
$ws.Range("M25").Value = 1.055336392457733
